# Daily attendance processing - 2025-11-17 18:29:23
# Update "Recorded By" attendance recorder lists (re-ordered / updated) and
# related derived attendance statistics on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ANATOMY, session 1) - reorder recorders
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# Row 3 (ANATOMY, session 2) - reorder recorders
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 (ANATOMY, session 3) - updated recorders list and new attendance count
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H4").Value = "73/251"

# Row 9 (HISTOLOGY, session 1) - reorder recorders
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 10 stats - updated average attendance percentage
$ws.Range("L10").Value = "25.5%"

# Row 15 (PARASITOLOGY, session 2) - reorder recorders
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# Row 15 stats - updated average attendance percentage
$ws.Range("S15").Value = "25.5%"

# Row 28 (PHYSIOLOGY, session 1) - reorder recorders
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
